# Append: 2025-11-26 01:50 JST
# The scraper re-ran; the previous top-of-list job posting (row 6,
# "縫製工場向けPL・CF可視化アプリ") is no longer present in the new
# scrape, so that row is removed and everything below it shifts up by
# one. All remaining rows get their "取得日時" (fetched-at) timestamp
# bumped to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove the row that dropped out of this scrape; Excel's row delete
# shifts rows 7:19 up to 6:18 and fixes up the hyperlink relationships
# automatically.
$ws.Rows.Item(6).Delete()

$newTimestamp = "2025-11-26 01:50:25"

# Update the fetch timestamp for every remaining data row (now rows 2-18).
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
